# Atualizacao de bases das ligas, do dia: 21-02-2024 as 23:25
#
# The underlying match records (columns B:AC, i.e. "id" through "PL_AhUnder")
# for several rows got re-matched against the correct fixtures. Column A
# (the sequential row index) stays put; only the data describing each
# fixture (id, teams, score, odds, ...) moves between rows.
#
# Each block below reads the current B:AC contents of the rows involved,
# then writes them back in the new arrangement:
#   - most blocks are simple 2-row swaps
#   - the block for rows 177/179/180 is a 3-way rotation (row 178 is
#     untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 32 <-> 33 ---
$rA = $ws.Range("B32:AC32")
$rB = $ws.Range("B33:AC33")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 40 <-> 41 ---
$rA = $ws.Range("B40:AC40")
$rB = $ws.Range("B41:AC41")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 67 <-> 68 ---
$rA = $ws.Range("B67:AC67")
$rB = $ws.Range("B68:AC68")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 77 <-> 78 ---
$rA = $ws.Range("B77:AC77")
$rB = $ws.Range("B78:AC78")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 137 <-> 138 ---
$rA = $ws.Range("B137:AC137")
$rB = $ws.Range("B138:AC138")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 144 <-> 145 ---
$rA = $ws.Range("B144:AC144")
$rB = $ws.Range("B145:AC145")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 168 <-> 169 ---
$rA = $ws.Range("B168:AC168")
$rB = $ws.Range("B169:AC169")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# --- rows 177, 179, 180 (3-way rotation; row 178 stays as-is) ---
# new(177) = old(180); new(179) = old(177); new(180) = old(179)
$r177 = $ws.Range("B177:AC177")
$r179 = $ws.Range("B179:AC179")
$r180 = $ws.Range("B180:AC180")
$v177 = $r177.Value2
$v179 = $r179.Value2
$v180 = $r180.Value2
$r177.Value2 = $v180
$r179.Value2 = $v177
$r180.Value2 = $v179
